# Automatische test-sync: 2025-07-29 21:55:50
# Append a new test-mail log entry (row 15) to the "Logs" sheet and
# refresh the "Dashboard" category summary to reflect it.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- New row 15 on the Logs sheet -----------------------------------------
$logs.Range("A15").Value = "Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Range("B15").Value = "mailmind.test@zohomail.eu"
$logs.Range("C15").Value = "Testmail #13: Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Range("D15").Value = "Productinformatie"

$antwoord = "Beste klant," + [char]10 + `
    "Bedankt voor je bericht. Helaas kan ik je momenteel niet de datasheet van de VentiQ-250 sturen, aangezien ik niet beschik over het specifieke materiaal waar je naar vraagt." + [char]10 + `
    "Ik raad je aan om contact op te nemen met onze verkoopafdeling of de klantenservice, zodat zij je verder kunnen helpen met het verkrijgen van de juiste informatie." + [char]10 + `
    "Met vriendelijke groet," + [char]10 + `
    "[Naam]  " + [char]10 + `
    "E-mailassistent bij [Bedrijfsnaam]"
$logs.Range("E15").Value = $antwoord

$logs.Range("F15").Value = "2025-07-29 21:55:30"
$logs.Range("G15").Value = "Ja"
$logs.Range("H15").Value = "Nee"
$logs.Range("I15").Value = "Ja"
$logs.Range("J15").Value = "Nee"

# Multi-line content in E15 makes the engine auto-pin an explicit row
# height; re-run autofit so the row goes back to a plain, unpinned height
# (matching how the rest of the sheet's rows -- which also hold multi-line
# "Antwoord" text -- are stored).
$logs.Rows(15).AutoFit()

# --- Expand conditional formatting ranges to include the new row ----------
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $logs.Range($col + "2:" + $col + "14")
    $newRange = $logs.Range($col + "2:" + $col + "15")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard: Productinformatie count goes 3 -> 4 and now outranks ------
# --- "Intern verzoek / Actie voor medewerker" (still 3), so the two rows --
# --- swap places in the sorted-by-count summary table. --------------------
$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 4
$dash.Range("A4").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B4").Value = 3
